# Updated presentation for Dr Herceg
# Resize/reposition the Caesar-cipher illustration picture on slide 2
# (the shift from the old box to the new, slightly smaller/centered box).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)   # "Picture 2"

# Target EMU values (from the authored edit):
#   off: x=1362811, y=2609813
#   ext: cx=3106447, cy=1862579
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points
# (1 pt = 12700 EMU); the literals below are chosen so the point -> EMU
# conversion lands exactly on the target EMU values.
$sh.Left   = 107.30794907
$sh.Top    = 205.49709321
$sh.Width  = 244.60211945
$sh.Height = 146.65976716
